# Update "想去人数" (want-to-go count) figures in both the "展览" sheet
# and the consolidated "全部类型" sheet to match the refreshed data pull.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 29
$ws1.Range("F5").Value = 307
$ws1.Range("F7").Value = 1054
$ws1.Range("F10").Value = 10
$ws1.Range("F13").Value = 13588
$ws1.Range("F15").Value = 21
$ws1.Range("F17").Value = 5584
$ws1.Range("F18").Value = 5591
$ws1.Range("F19").Value = 67

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 29
$ws4.Range("F21").Value = 307
$ws4.Range("F29").Value = 1054
$ws4.Range("F32").Value = 10
$ws4.Range("F35").Value = 13588
$ws4.Range("F37").Value = 21
$ws4.Range("F40").Value = 5584
$ws4.Range("F41").Value = 5591
$ws4.Range("F42").Value = 67
